# Normalize the "Recorded By" (column G) entries on the "Session Analysis
# Results" sheet: each cell holds a comma-separated list of names/emails
# that should be presented in (case-sensitive / ordinal) ascending order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column G is the 7th column ("Recorded By").
$col = 7

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $text = $cell.Value2

    if ($null -eq $text) { continue }
    if ($text -isnot [string]) { continue }
    if ($text -notmatch ',') { continue }

    $parts = @($text -split ',' | ForEach-Object { $_.Trim() })

    $list = New-Object System.Collections.Generic.List[string]
    foreach ($p in $parts) { [void]$list.Add($p) }
    $list.Sort([System.StringComparer]::Ordinal)

    $newText = ($list -join ', ')

    if ($newText -ne $text) {
        $cell.Value2 = $newText
    }
}
